# Applies the "Making changes to the test-cases" commit to
# an_Datasheet1.xlsx: two new values are entered in column Y (a new
# header "user_meghana" in Y1, and a record id "005q0000003GGfP" in Y2)
# and the sheet's viewport/selection is moved from T6 to Y2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (row 1) and new data cell (row 2) in column Y.
# These land as new shared-string entries (index 51 / 52) exactly as
# the target workbook does, and pick up the existing style (s="10")
# that was already on Y1/Y2.
$ws.Range("Y1").Value = "user_meghana"
$ws.Range("Y2").Value = "005q0000003GGfP"

# Move the window scroll position / active selection like the recorded
# session did (topLeftCell moved from P1 to T1, selection moved from
# T6 to Y2).
$excel.ActiveWindow.ScrollColumn = 20
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Y2").Select()
